$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raluca")

$ws.Range("C14").Value = "https://amzn.eu/d/fxzOrlU"
$ws.Range("C15").Value = "https://amzn.eu/d/hKdAdfZ"

$ws.Range("A14").Value = "Shiatsu Heat Massager"
$ws.Range("A15").Value = "Laneige Lip Sleeping Mask"

$ws.Range("B15").Value = "https://m.media-amazon.com/images/I/51R2w6yMcaL._SX450_.jpg"
$ws.Range("B14").Value = "https://m.media-amazon.com/images/I/61Hx-AZNQ-L._AC_SX450_.jpg"

$ws.Range("D14").Value = "39.97 EUR"
$ws.Range("D15").Value = "18.02 EUR"

[void]$ws.Range("D16").Select()
